# Sign-up automation workbook: replace the single "login" sheet with four
# sheets (invalidEmail, invalidPassword1/2/3) holding the new test-data
# fixtures, per the commit "Added 4 worksheets for Sign up Page".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the existing sheet as "invalidEmail" and wipe its old
#    login-fixture content (values, styles, hyperlinks) so we can lay
#    down the new data cleanly.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "invalidEmail"
$ws1.Hyperlinks.Delete()
$ws1.Cells.ClearContents()
$ws1.Cells.ClearFormats()

# ---------------------------------------------------------------------
# 2. Create the three additional sheets, in tab order, right after the
#    first one.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "invalidPassword1"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "invalidPassword2"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "invalidPassword3"

# ---------------------------------------------------------------------
# 3. invalidEmail sheet data
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "USERNAME"
$ws1.Range("B1").Value = "ERROR MESSAGE"

$ws1.Range("A2").Value = "t"
$ws1.Range("B2").Value = "Invalid Email Address"

$ws1.Range("A3").Value = "test1@"
$ws1.Range("B3").Value = "Invalid Email Address"

$ws1.Range("A4").Value = "test2@.com"
$ws1.Range("B4").Value = "Invalid Email Address"

$ws1.Range("A5").Value = "test3@gmail.c"
$ws1.Range("B5").Value = "Invalid Email Address"

$ws1.Range("A2").Style = "Hyperlink"
$ws1.Range("A3").Style = "Hyperlink"
$ws1.Range("A4").Style = "Hyperlink"
$ws1.Range("A5").Style = "Hyperlink"

# NOTE: `Hyperlinks.Add(..., TextToDisplay)` overwrites the backing cell's
# text with the display string (and a display-less Add on a multi-cell
# range stamps the anchor cell too) -- re-assert our intended cell text
# straight after each Add call so the stored <v> stays correct while the
# hyperlink's `display` attribute still comes out as requested.
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:test1@gmail.com", "", "", "test1@gmail.com")
$ws1.Range("A2").Value = "t"

$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:test1@gmail.com")
$ws1.Range("A3").Value = "test1@"

$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:test3@gmail.com")
$ws1.Range("A5").Value = "test3@gmail.c"

$ws1.Hyperlinks.Add($ws1.Range("A3:A5"), "mailto:test1@gmail.com", "", "", "test1@gmail.com")
$ws1.Range("A3").Value = "test1@"
$ws1.Range("A4").Value = "test2@.com"
$ws1.Range("A5").Value = "test3@gmail.c"

$ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:test2@gmail.com")
$ws1.Range("A4").Value = "test2@.com"

$ws1.Range("A2").Style = "Hyperlink"
$ws1.Range("A3").Style = "Hyperlink"
$ws1.Range("A4").Style = "Hyperlink"
$ws1.Range("A5").Style = "Hyperlink"

$ws1.Columns.Item(1).ColumnWidth = 22.666666666666668
$ws1.Columns.Item(2).ColumnWidth = 32.333333333333332

$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. invalidPassword1 sheet data
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "PASSWORD"
$ws2.Range("B1").Value = "ERROR MESSAGE"

$ws2.Range("A2").Value = "a"
$ws2.Range("B2").Value = "Must be 8-20 characters."

$ws2.Range("A3").Value = "abc1234"
$ws2.Range("B3").Value = "Must be 8-20 characters."

$ws2.Range("A4").Value = "Abc45678901234567890"
$ws2.Range("B4").Value = "Must be 8-20 characters."

$ws2.Columns.Item(1).ColumnWidth = 22.0
$ws2.Columns.Item(2).ColumnWidth = 21.833333333333332

$ws2.Range("B9").Select()

# ---------------------------------------------------------------------
# 5. invalidPassword2 sheet data
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "PASSWORD"
$ws3.Range("B1").Value = "ERROR MESSAGE"

$ws3.Range("A2").Value = "Abcdefgh"
$ws3.Range("B2").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws3.Range("A3").Value = 12345678
$ws3.Range("B3").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws3.Range("A4").Value = "!@#$%^&*"
$ws3.Range("B4").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws3.Range("A5").Value = "1234567A"
$ws3.Range("B5").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws3.Range("A6").Value = "ABC1234$"
$ws3.Range("B6").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws3.Range("A4").Style = "Hyperlink"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "mailto:test1@gmail.com")
$ws3.Range("A4").Style = "Hyperlink"

$ws3.Columns.Item(1).ColumnWidth = 25.333333333333332
$ws3.Columns.Item(2).ColumnWidth = 70.833333333333329

$ws3.Range("B19").Select()

# ---------------------------------------------------------------------
# 6. invalidPassword3 sheet data
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "PASSWORD"
$ws4.Range("B1").Value = "ERROR MESSAGE"

$ws4.Range("A2").Value = "A"
$ws4.Range("B2").Value = "Password Required."

$ws4.Range("A3").Value = "Abc12345"
$ws4.Range("B3").Value = "Must be 8-20 characters."

$ws4.Range("A4").Value = "ABC1234$"
$ws4.Range("B4").Value = "Must contain one upper & lower case letter and a non-letter (number or symbol.)"

$ws4.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws4.Columns.Item(2).ColumnWidth = 17.833333333333332

$ws4.Range("B4").Select()

# ---------------------------------------------------------------------
# 7. Final selection/active sheet: invalidEmail is the tab shown when
#    the workbook is opened (tabSelected="1" in the saved file).
# ---------------------------------------------------------------------
$ws1.Range("E12").Select()
$ws1.Activate()
